# Generate Report for Handoff
# Inserts a new localization-status row for file
# "5c528238-8f6e-4fd6-bebd-de0b1ce0df4b.md" into every sheet (Overview,
# zh-cn, de-de), positioned between the "d529a662-..." and "64ae4af2-..."
# entries, and rebuilds the hyperlinks on each sheet to match.

$wb = $excel.ActiveWorkbook

$newMdName  = "5c528238-8f6e-4fd6-bebd-de0b1ce0df4b.md"
$newBase    = "5c528238-8f6e-4fd6-bebd-de0b1ce0df4b"
$newXlfHash = "4cc5096ef129e14167d6a7c2a4fcaf7db018d070"
$newStatus  = "Ready for handoff"

# Row index (1-based) at which the new record is inserted on every sheet.
$insertRow = 7

# ======================================================================
# Sheet 1: "Overview"  (columns: File Name | zh-cn | de-de | Latest Handoff Date)
# ======================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item($insertRow).Insert() | Out-Null

$ws.Range("A$insertRow").Value = $newMdName
$ws.Range("B$insertRow").Value = $newStatus
$ws.Range("C$insertRow").Value = $newStatus
$ws.Range("D$insertRow").Value = "2016-03-25 02:53:02"

# Rebuild hyperlinks top to bottom (row-insert does not shift existing
# Hyperlink ranges in this engine, so clear + recreate from scratch).
$ws.Cells.Hyperlinks.Delete()

$overviewLinks = @(
    @{ Row = 2; File = "bd1e060a-d5db-48a9-95cf-5645fc0d341e.md" },
    @{ Row = 3; File = "0f4cbed8-f610-4895-b315-31b06abe215a.md" },
    @{ Row = 4; File = "1a7f879d-57a5-46fa-b42e-15137a4100ba.md" },
    @{ Row = 5; File = "7e2d49fa-7c1e-4196-990c-84cc4565ed82.md" },
    @{ Row = 6; File = "d529a662-e74c-420b-b74a-f561886915c6.md" },
    @{ Row = 7; File = $newMdName },
    @{ Row = 8; File = "64ae4af2-4a30-4e99-91db-fb8621460b73.md" },
    @{ Row = 9; File = "9d500b84-b057-42fb-88e9-5ce55a2109ec.md" }
)
foreach ($link in $overviewLinks) {
    $addr = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$($link.File)"
    $ws.Hyperlinks.Add($ws.Range("A$($link.Row)"), $addr, [Type]::Missing, [Type]::Missing, $link.File) | Out-Null
}

# ======================================================================
# Sheet 2: "zh-cn"
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime |
#          F Latest Target File | G Latest Handback File |
#          H Latest Handback DateTime | I Reference Tokens |
#          J Handoff Reason | K Dependency From | L Error Detail
# ======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$lang = "zh-cn"

$wsZh.Rows.Item($insertRow).Insert() | Out-Null

$wsZh.Range("A$insertRow").Value = $newMdName
$wsZh.Range("B$insertRow").Value = ".md"
$wsZh.Range("C$insertRow").Value = $newStatus
$wsZh.Range("D$insertRow").Value = "$newBase.$newXlfHash.$lang.xlf"
$wsZh.Range("E$insertRow").Value = "2016-03-25 02:52:57"
$wsZh.Range("H$insertRow").Value = "0001-01-01 00:00:00"
$wsZh.Range("J$insertRow").Value = "Include"

$wsZh.Cells.Hyperlinks.Delete()

$zhRows = @(
    @{ Row = 2; File = "bd1e060a-d5db-48a9-95cf-5645fc0d341e"; Hash = "7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5"; WithTarget = $true  },
    @{ Row = 3; File = "0f4cbed8-f610-4895-b315-31b06abe215a"; Hash = "531d83604c9be95f6aa0c59bcd283eca044880ba"; WithTarget = $false },
    @{ Row = 4; File = "1a7f879d-57a5-46fa-b42e-15137a4100ba"; Hash = "4d45a5878cdf42b5f5292780578593e11129bf76"; WithTarget = $false },
    @{ Row = 5; File = "7e2d49fa-7c1e-4196-990c-84cc4565ed82"; Hash = "bcb225122fa9e924d07655517791494c86c1b676"; WithTarget = $true  },
    @{ Row = 6; File = "d529a662-e74c-420b-b74a-f561886915c6"; Hash = "c0f2756b5b323d58b52b7d0cb2af06df35ff9e4c"; WithTarget = $false },
    @{ Row = 7; File = $newBase;                                Hash = $newXlfHash;                                         WithTarget = $false },
    @{ Row = 8; File = "64ae4af2-4a30-4e99-91db-fb8621460b73"; Hash = "89d73eaeb163482b25d39c39b4a8923751f85584"; WithTarget = $false },
    @{ Row = 9; File = "9d500b84-b057-42fb-88e9-5ce55a2109ec"; Hash = "0b13aaf1d399248bea0c200d555b3d38fe629713"; WithTarget = $false }
)

foreach ($r in $zhRows) {
    $mdName  = "$($r.File).md"
    $xlfName = "$($r.File).$($r.Hash).$lang.xlf"

    $mdAddr  = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName"
    $xlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/$xlfName"

    $wsZh.Hyperlinks.Add($wsZh.Range("A$($r.Row)"), $mdAddr, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $wsZh.Hyperlinks.Add($wsZh.Range("D$($r.Row)"), $xlfAddr, [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null

    if ($r.WithTarget) {
        $wsZh.Hyperlinks.Add($wsZh.Range("F$($r.Row)"), $mdAddr, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
        $wsZh.Hyperlinks.Add($wsZh.Range("G$($r.Row)"), $xlfAddr, [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null
    }
}

# ======================================================================
# Sheet 3: "de-de"  (same layout as "zh-cn")
# ======================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$lang = "de-de"

$wsDe.Rows.Item($insertRow).Insert() | Out-Null

$wsDe.Range("A$insertRow").Value = $newMdName
$wsDe.Range("B$insertRow").Value = ".md"
$wsDe.Range("C$insertRow").Value = $newStatus
$wsDe.Range("D$insertRow").Value = "$newBase.$newXlfHash.$lang.xlf"
$wsDe.Range("E$insertRow").Value = "2016-03-25 02:53:02"
$wsDe.Range("H$insertRow").Value = "0001-01-01 00:00:00"
$wsDe.Range("J$insertRow").Value = "Include"

$wsDe.Cells.Hyperlinks.Delete()

$deRows = @(
    @{ Row = 2; File = "bd1e060a-d5db-48a9-95cf-5645fc0d341e"; Hash = "7880a49fe502cc1f7a2ce60e119dd66b8f1e69b5"; WithTarget = $true  },
    @{ Row = 3; File = "0f4cbed8-f610-4895-b315-31b06abe215a"; Hash = "531d83604c9be95f6aa0c59bcd283eca044880ba"; WithTarget = $false },
    @{ Row = 4; File = "1a7f879d-57a5-46fa-b42e-15137a4100ba"; Hash = "4d45a5878cdf42b5f5292780578593e11129bf76"; WithTarget = $false },
    @{ Row = 5; File = "7e2d49fa-7c1e-4196-990c-84cc4565ed82"; Hash = "bcb225122fa9e924d07655517791494c86c1b676"; WithTarget = $true  },
    @{ Row = 6; File = "d529a662-e74c-420b-b74a-f561886915c6"; Hash = "c0f2756b5b323d58b52b7d0cb2af06df35ff9e4c"; WithTarget = $false },
    @{ Row = 7; File = $newBase;                                Hash = $newXlfHash;                                         WithTarget = $false },
    @{ Row = 8; File = "64ae4af2-4a30-4e99-91db-fb8621460b73"; Hash = "89d73eaeb163482b25d39c39b4a8923751f85584"; WithTarget = $false },
    @{ Row = 9; File = "9d500b84-b057-42fb-88e9-5ce55a2109ec"; Hash = "0b13aaf1d399248bea0c200d555b3d38fe629713"; WithTarget = $false }
)

foreach ($r in $deRows) {
    $mdName  = "$($r.File).md"
    $xlfName = "$($r.File).$($r.Hash).$lang.xlf"

    $mdAddr  = "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$mdName"
    $xlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/$xlfName"

    $wsDe.Hyperlinks.Add($wsDe.Range("A$($r.Row)"), $mdAddr, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
    $wsDe.Hyperlinks.Add($wsDe.Range("D$($r.Row)"), $xlfAddr, [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null

    if ($r.WithTarget) {
        $wsDe.Hyperlinks.Add($wsDe.Range("F$($r.Row)"), $mdAddr, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
        $wsDe.Hyperlinks.Add($wsDe.Range("G$($r.Row)"), $xlfAddr, [Type]::Missing, [Type]::Missing, $xlfName) | Out-Null
    }
}
